# Applies the Swahili (Kenya) translation edits described by the diff.
$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Exact "Video Title" "Kichwa cha Video"
Replace-Exact "Topic" "Mada"
Replace-Exact "Aim(s)" "Malengo"
Replace-Exact "Length" "Urefu"
Replace-Exact "Camp Location" "Mahali pa Kambi"
Replace-Exact "Facilitators" "Wawezeshaji"
Replace-Exact "N. of students" "N. ya wanafunzi"
Replace-Exact "Date" "Tarehe"
Replace-Exact "Resources" "Rasilimali"
Replace-Exact "needed" "inahitajika"
Replace-Exact "Preparations" "Maandalizi"
Replace-Exact "Video time" "Muda wa video"
Replace-Exact "What facilitator does" "Mwezeshaji anafanya nini"
Replace-Exact "What learners do" "Wanachofanya wanafunzi"
Replace-Exact "General VMC Video Introduction" "Utangulizi Mkuu wa Video ya VMC"
Replace-Exact "Introduction of the first experiment" "Utangulizi wa jaribio la kwanza"
Replace-Exact "Assist the process, provoke thoughts" "Kusaidia mchakato, kuchochea mawazo"
